$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to be treated as text so numeric-looking
# values (e.g. "13.70", "1.00") keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.850.46"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.212.25"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "291.58"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "86.81"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.467"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "30.48"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "0.0777"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "49.99"
$ws.Range("E12").Value = "  +5.92%  "
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "2.559.74"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "2.268.49"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "13.70"
$ws.Range("D18").Value = "0.729"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "39.792.54"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "11.23"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "5.74"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "65.50"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "236.36"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "1.83"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "23.40"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("D30").Value = "9.21"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("D32").Value = "31.75"
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "0.0709"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("D36").Value = "2.90"
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").Value = "0.0980"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").Value = "1.70"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "15.17"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("D42").Value = "2.107.40"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").Value = "0.0268"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.09"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "17.67"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "9.85"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "2.70"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("D49").Value = "2.433.98"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "1.48"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("D51").Value = "88.37"
$ws.Range("E51").Value = "  -0.93%  "
